# "Ajout du system de drag and drop"
#
# The journal table (Tableau1, on the "Activités" sheet) gained a new
# drag-and-drop entry on row 29: a Date ("Début" day) and a start Time
# were dropped into the row, which feeds the existing "Temps" calculated
# column (Fin blank -> NOW() - Début). The view was also left scrolled to
# where that new row now sits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: the newly dropped activity - Date (A29) + Début time (B29).
# 44279 -> 24 March 2021 ; 0.36805555555555558 -> 08:50:00.
$ws.Range("A29").Value = 44279
$ws.Range("B29").Value = 0.36805555555555558

# The "Temps" column (D) keeps its live calculated-column formula; it
# recalculates automatically now that Début is filled in (Fin is still
# blank, so it reports NOW() - Début).

# Move the viewport/selection to follow the row that was just dropped in.
[void]$ws.Range("E29").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
